# Updated cryptos list values (price + 1h volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number: force text format first so
# Excel does not silently convert them to the numeric type (source data keeps
# these as text cells, e.g. "306.76").
$textCells = @("D5", "D6", "D7", "D9", "D10", "D13", "D17", "D19", "D22", "D23", "D25", "D27", "D28", "D29", "D31", "D32", "D35", "D36", "D37", "D38", "D41", "D42", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values
$ws.Range("D2").Value = "42.364.82"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.272.67"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "306.76"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "97.40"
$ws.Range("E6").Value = "  +5.15%  "
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  +9.09%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "6.67"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "2.624.15"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "2.251.05"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "42.249.20"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "67.57"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "240.42"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "1.95"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "23.83"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "37.33"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "159.85"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +4.79%  "
$ws.Range("D35").Value = "0.0741"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "17.03"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  +4.42%  "
$ws.Range("D42").Value = "2.44"
$ws.Range("E42").Value = "  +14.30%  "
$ws.Range("D43").Value = "1.998.36"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "18.89"
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "9.98"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").Value = "52.96"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "72.16"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "91.54"
$ws.Range("E51").Value = "  +0.53%  "

# Drop the temporary number format again so styling matches the source (General)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
